$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1 & 3: "英語" -> "英语" (replaces both the hyperlink text and the plain paragraph text,
#         Word's Find.Execute with ReplaceAll will hit every occurrence)
Replace-Text "英語" "英语"

# 2: language list with traditional -> simplified characters
# (match starts at "/" rather than the leading space, because the leading space sits
#  exactly at the hyperlink/run boundary and Find would otherwise mis-attribute the
#  replacement's formatting to the preceding hyperlink run instead of this run's ff0000)
Replace-Text "/ 葡萄牙語 / 法語 / 泰語 / 越南語 / 西班牙語" "/ 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语"

# 4: 簡介 -> 简介
Replace-Text "簡介" "简介"

# 5: description sentence (also resolves the only other "。 " occurrence)
Replace-Text "發送給目標國家中那些文件未通過我們驗證流程的合作夥伴的電子郵件。 將通過 customer.io 發送" "一封发送给目标国家中未通过我们验证流程的合作伙伴的电子邮件。 将通过 customer.io 发送"

# 6: 目標受眾 -> 目标受众
Replace-Text "目標受眾" "目标受众"

# 7: 提交錯誤/不完整文檔的被邀請合作夥伴 -> 提交了错误/不完整文件的邀请合作伙伴
Replace-Text "提交錯誤/不完整文檔的被邀請合作夥伴" "提交了错误/不完整文件的邀请合作伙伴"

# 8: 主題行 -> 主题行
Replace-Text "主題行" "主题行"

# 9: [事件名稱] -> [事件名称]
Replace-Text "[事件名稱]" "[事件名称]"

# 10: " — 文件驗證失敗 " -> " — 文档验证失败 "
Replace-Text " — 文件驗證失敗 " " — 文档验证失败 "

# 11: 啊哦！ 文檔無法驗證 -> 啊哦！ 文件无法验证
Replace-Text "啊哦！ 文檔無法驗證" "啊哦！ 文件无法验证"

# 12: [合作夥伴姓名] -> [合作伙伴姓名]
Replace-Text "[合作夥伴姓名]" "[合作伙伴姓名]"

# 13: regret sentence -> English replacement text
Replace-Text "很遺憾地通知您，您的文檔未通過驗證流程，因為我們發現以下問題： " "We regret to inform you that your documents have failed our verification process as we found the following issues with them: "

# 14: 您的疫苗接種證明副本 -> 疫苗接种证书副本
Replace-Text "您的疫苗接種證明副本" "疫苗接种证书副本"

# 15: : 文檔不清楚 -> : 文件不清楚
Replace-Text ": 文檔不清楚" ": 文件不清楚"

# 16: [文檔 2] -> [文件 2]
Replace-Text "[文檔 2]" "[文件 2]"

# 17: : [問題] -> : [problem]
Replace-Text ": [問題]" ": [problem]"

# 18: 請在  -> 请在 
Replace-Text "請在 " "请在 "

# 19: resubmit sentence
Replace-Text " 之前重新提交上述文檔，以便我們進行必要的安排。" " 之前重新提交上述文件，以便我们进行必要的安排。"

# 20: 如有任何疑問，請通過  -> 如有任何疑问，请通过 
Replace-Text "如有任何疑問，請通過 " "如有任何疑问，请通过 "

# 21: [電子郵件地址] -> [电子邮件地址]
Replace-Text "[電子郵件地址]" "[电子邮件地址]"

# 22: [WHATSAPP 號碼] -> [WHATSAPP 号码]
Replace-Text "[WHATSAPP 號碼]" "[WHATSAPP 号码]"

# 23: contact regional manager sentence (comma after 经理 removed)
Replace-Text " (WhatsApp) 聯繫您的區域經理, " " (WhatsApp) 联系您的区域经理 "

# 24: [姓名] -> [NAME] (this run's entire content, so the run split is preserved)
Replace-Text "[姓名]" "[NAME]"

# 25: the trailing "。 " right after [NAME] becomes " 。 " (leading space added).
# A plain ReplaceAll on "。 " would be ambiguous since the translated sentence in step 5
# also contains "。 " (一封...电子邮件。 将通过...). To target only the final occurrence,
# locate [NAME] first and then scope the Find to the small range that immediately follows it,
# using ReplaceOne (not ReplaceAll) with wrap disabled so the whole document isn't rescanned.
$nameRange = $d.Content
$nameRange.Find.Execute("[NAME]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tailRange = $d.Range($nameRange.End, $nameRange.End + 5)
$tailRange.Find.Execute("。 ", $true, $false, $false, $false, $false, $true, 0, $false, " 。 ", 1) | Out-Null
